$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1115358.5
$ws.Range("I6").Value = 2000045.6
$ws.Range("J6").Value = 9499.5
$ws.Range("K6").Value = 6000136.800000001
$ws.Range("L6").Value = 28498.5
$ws.Range("M6").Value = -6000024.800000001
$ws.Range("N6").Value = -28722.5
$ws.Range("H62").Value = 54976.19
$ws.Range("I62").Value = 83353.84
$ws.Range("K62").Value = 83353.84
$ws.Range("M62").Value = -82729.84
$ws.Range("H65").Value = 54976.19
$ws.Range("I65").Value = 83353.84
$ws.Range("K65").Value = 416769.2
$ws.Range("M65").Value = -413649.2
$ws.Range("H141").Value = 1362
$ws.Range("I141").Value = 1362
$ws.Range("K141").Value = 4086
$ws.Range("M141").Value = 1094

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2385
$ws.Range("I132").Value = 2135.9524
$ws.Range("K132").Value = 6407.8572
$ws.Range("M132").Value = -3877.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31148
$ws.Range("H99").Value = 2044.2222
$ws.Range("I99").Value = 991.6667
$ws.Range("J99").Value = 4149.3335
$ws.Range("K99").Value = 991.6667
$ws.Range("L99").Value = 4149.3335
$ws.Range("M99").Value = 506.3333
$ws.Range("N99").Value = -7145.3335
$ws.Range("H107").Value = 1679.6875
$ws.Range("J107").Value = 1909.2
$ws.Range("L107").Value = 1909.2
$ws.Range("N107").Value = -5749.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 57810.105
$ws.Range("I31").Value = 73864.07000000001
$ws.Range("J31").Value = 12859
$ws.Range("K31").Value = 73864.07000000001
$ws.Range("L31").Value = 12859
$ws.Range("M31").Value = -73569.07000000001
$ws.Range("N31").Value = -13449
$ws.Range("H32").Value = 16754.25
$ws.Range("I32").Value = 16754.25
$ws.Range("K32").Value = 16754.25
$ws.Range("M32").Value = -16438.25
$ws.Range("H34").Value = 57810.105
$ws.Range("I34").Value = 73864.07000000001
$ws.Range("J34").Value = 12859
$ws.Range("K34").Value = 73864.07000000001
$ws.Range("L34").Value = 12859
$ws.Range("M34").Value = -73662.07000000001
$ws.Range("N34").Value = -13263
$ws.Range("H106").Value = 34223.332
$ws.Range("J106").Value = 34223.332
$ws.Range("L106").Value = 34223.332
$ws.Range("N106").Value = -36747.332
$ws.Range("H107").Value = 815.2083
$ws.Range("I107").Value = 884.58826
$ws.Range("K107").Value = 884.58826
$ws.Range("M107").Value = 1035.41174
$ws.Range("H122").Value = 1487.6666
$ws.Range("I122").Value = 1481.5
$ws.Range("K122").Value = 4444.5
$ws.Range("M122").Value = -1994.5
$ws.Range("H132").Value = 2670.6843
$ws.Range("I132").Value = 2553.6758
$ws.Range("K132").Value = 7661.0274
$ws.Range("M132").Value = -5131.0274
$ws.Range("H134").Value = 19360.258
$ws.Range("I134").Value = 8487.23
$ws.Range("J134").Value = 75900
$ws.Range("K134").Value = 25461.69
$ws.Range("L134").Value = 227700
$ws.Range("M134").Value = -22926.69
$ws.Range("N134").Value = -232770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23029082
$ws.Range("I4").Value = 26804492
$ws.Range("J4").Value = 2893555.5
$ws.Range("K4").Value = 80413476
$ws.Range("L4").Value = 8680666.5
$ws.Range("M4").Value = -80413364
$ws.Range("N4").Value = -8680890.5
$ws.Range("H23").Value = 91.09090999999999
$ws.Range("I23").Value = 102.5
$ws.Range("K23").Value = 307.5
$ws.Range("M23").Value = -72.5
$ws.Range("H113").Value = 519.5
$ws.Range("I113").Value = 469.8
$ws.Range("J113").Value = 555
$ws.Range("K113").Value = 1409.4
$ws.Range("L113").Value = 1665
$ws.Range("M113").Value = 760.5999999999999
$ws.Range("N113").Value = -6005
$ws.Range("H122").Value = 2675.9285
$ws.Range("I122").Value = 2615.5715
$ws.Range("J122").Value = 2736.2856
$ws.Range("K122").Value = 23540.1435
$ws.Range("L122").Value = 24626.5704
$ws.Range("M122").Value = -21090.1435
$ws.Range("N122").Value = -29526.5704
$ws.Range("H124").Value = 27141.857
$ws.Range("J124").Value = 27141.857
$ws.Range("L124").Value = 81425.571
$ws.Range("N124").Value = -91245.571
$ws.Range("H129").Value = 542.0769
$ws.Range("I129").Value = 542.0769
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1626.2307
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3373.7693
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 65768.88
$ws.Range("J131").Value = 9996.066000000001
$ws.Range("L131").Value = 29988.198
$ws.Range("N131").Value = -40068.198
$ws.Range("H132").Value = 1340.4814
$ws.Range("I132").Value = 1138.826
$ws.Range("K132").Value = 10249.434
$ws.Range("M132").Value = -7719.434000000001
$ws.Range("H137").Value = 4109.5
$ws.Range("I137").Value = 2715.6365
$ws.Range("K137").Value = 8146.9095
$ws.Range("M137").Value = -3046.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5739.483
$ws.Range("I70").Value = 4676.885
$ws.Range("K70").Value = 4676.885
$ws.Range("M70").Value = -4406.885
$ws.Range("H73").Value = 5739.483
$ws.Range("I73").Value = 4676.885
$ws.Range("K73").Value = 4676.885
$ws.Range("M73").Value = -3740.885
$ws.Range("H122").Value = 2345.125
$ws.Range("I122").Value = 2318.4167
$ws.Range("K122").Value = 6955.250100000001
$ws.Range("M122").Value = -4505.250100000001
$ws.Range("H132").Value = 266636.44
$ws.Range("I132").Value = 266636.44
$ws.Range("K132").Value = 799909.3200000001
$ws.Range("M132").Value = -797379.3200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 230
$ws.Range("I13").Value = 230
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 230
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -90
$ws.Range("N13").ClearContents()
$ws.Range("H22").Value = 1152.2
$ws.Range("J22").Value = 1346
$ws.Range("L22").Value = 1346
$ws.Range("N22").Value = -1936
$ws.Range("H27").Value = 1152.2
$ws.Range("J27").Value = 1346
$ws.Range("L27").Value = 1346
$ws.Range("N27").Value = -1560
$ws.Range("H40").Value = 6036.5713
$ws.Range("I40").Value = 5651.2
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 5651.2
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -5515.2
$ws.Range("N40").Value = -7272
$ws.Range("H55").Value = 146.04762
$ws.Range("I55").Value = 143
$ws.Range("J55").Value = 175
$ws.Range("K55").Value = 143
$ws.Range("L55").Value = 175
$ws.Range("M55").Value = 30
$ws.Range("N55").Value = -521
$ws.Range("H122").Value = 10976.846
$ws.Range("I122").Value = 7355.5557
$ws.Range("J122").Value = 19124.75
$ws.Range("K122").Value = 22066.6671
$ws.Range("L122").Value = 57374.25
$ws.Range("M122").Value = -19616.6671
$ws.Range("N122").Value = -62274.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3886.2942
$ws.Range("I122").Value = 3267.6
$ws.Range("J122").Value = 4770.143
$ws.Range("K122").Value = 9802.799999999999
$ws.Range("L122").Value = 14310.429
$ws.Range("M122").Value = -7352.799999999999
$ws.Range("N122").Value = -19210.429
$ws.Range("H132").Value = 5040.5757
$ws.Range("I132").Value = 3938.3333
$ws.Range("J132").Value = 10000.667
$ws.Range("K132").Value = 11814.9999
$ws.Range("L132").Value = 30002.001
$ws.Range("M132").Value = -9284.999899999999
$ws.Range("N132").Value = -35062.001
$ws.Range("H136").Value = 1892.3611
$ws.Range("I136").Value = 1682.0322
$ws.Range("K136").Value = 5046.096600000001
$ws.Range("M136").Value = -2496.096600000001

Write-Output "Applied Coeurl_Profits updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
